# Applies updated bus voltage magnitude (vm_pu) results for the 380 kV case.
# Columns B,C,D,E,F,I,J,K,L,M,N for rows 2-25 are updated with new values;
# the first column (A, the index) and columns G/H are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters (in the same order the values appear in each data row below)
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# Each entry: row number followed by the 11 new values for columns B..N (skipping G/H positions not listed)
$data = @(
    "2,1.02,1.022545013070072,1.028534243023331,0.9926147277508489,1.020994917831002,1.032211687326401,1.027730249084565,1.031350780310421,0.9955398523336033,1.023833557202326,1.01332284287209",
    "3,1.02,1.02334415576611,1.029136085777268,0.9936372048519304,1.02242941244794,1.032394726607857,1.02816826774634,1.031761134877081,0.9963617723202692,1.025072689561317,1.013467209998588",
    "4,1.02,1.023861143223613,1.029525282048165,0.9942998659930995,1.023357568558902,1.032511655339756,1.028450929776744,1.032025731677615,0.9968940712668345,1.025873930793942,1.013560366477664",
    "5,1.02,1.024078457364732,1.029688842548306,0.9945786998346017,1.023747755107751,1.032560450303047,1.028569577245361,1.032136744911038,0.997117960005301,1.026210641441592,1.013599467367713",
    "6,1.02,1.024114943737732,1.029716301666686,0.9946255319796338,1.023813268646867,1.032568621959934,1.028589487895653,1.032155371411488,0.9971555583673453,1.026267169035239,1.013606028933334",
    "7,1.02,1.023864047094348,1.029527467779295,0.9943035907982488,1.023362782286622,1.03251230876287,1.028452515872188,1.032027215919913,0.9968970624462087,1.02587843044981,1.01356088918986",
    "8,1.02,1.022815109689144,1.028737686880995,0.9929600610674301,1.021479726123447,1.032273858710341,1.027878437754653,1.031489653766051,0.995817528259106,1.024252445574596,1.013371685889979",
    "9,1.02,1.020965927602569,1.027344228895634,0.9906006454969559,1.018160937904734,1.031842132285142,1.026860995573392,1.030535303885412,0.9939188001724441,1.021382828666245,1.013036310991647",
    "10,1.02,1.019732629087644,1.026414132234902,0.989033133672735,1.015947781958512,1.031546573382374,1.026178795240223,1.029894333275141,0.9926553831429383,1.019466588640938,1.012811408893734",
    "11,1.02,1.019198484115299,1.026011135714611,0.988355674866747,1.014989260749221,1.031416761340392,1.025882472954668,1.029615668995271,0.9921088820399291,1.018636043144438,1.012713712657629",
    "12,1.02,1.019000061637497,1.025861406968857,0.9881042295826724,1.014633187279739,1.031368268265105,1.02577226696949,1.029511992920068,0.9919059725120875,1.018327417495245,1.012677377142",
    "13,1.02,1.019042624704285,1.025893525992812,0.9881581567098651,1.014709567949036,1.031378682638618,1.025795912799171,1.029534239373735,0.9919494934313052,1.018393624419399,1.012685173348722",
    "14,1.02,1.019182082800097,1.025998759857382,0.9883348863814464,1.014959828365643,1.031412758497971,1.025873366121272,1.029607102517938,0.9920921077337197,1.018610534604723,1.012710710106301",
    "15,1.02,1.019268005266251,1.026063592921222,0.9884438009545853,1.015114017162272,1.031433717314588,1.025921069299259,1.029651973691046,0.9921799884222134,1.018744163575938,1.012726437949508",
    "16,1.02,1.019768076080182,1.026440872420815,0.9890781214508737,1.016011391048568,1.031555149984518,1.026198441709108,1.029912803754332,0.9926916645766087,1.019521692016603,1.01281788610431",
    "17,1.02,1.020081725942992,1.026677461339093,0.989476357848556,1.016574230597981,1.031630830931981,1.026372182635647,1.030076116106948,0.9930127773699352,1.020009197999633,1.012875165626194",
    "18,1.02,1.020264661123146,1.026815434645199,0.9897087662937556,1.016902505414467,1.031674797462171,1.026473433507174,1.030171265365293,0.9932001317071769,1.020293474966158,1.012908545677474",
    "19,1.02,1.020327035327973,1.026862475706444,0.9897880325774034,1.017014435391409,1.031689758900722,1.026507942287139,1.03020369045761,0.9932640239640975,1.020390393114163,1.012919922299208",
    "20,1.02,1.020048075464159,1.02665208016692,0.9894336180360679,1.016513845338247,1.031622729376513,1.026353551110147,1.030058605412142,0.9929783193494215,1.019956901215049,1.01286902319069",
    "21,1.02,1.019141016338543,1.025967772145679,0.9882828385668249,1.014886133940531,1.031402731598527,1.025850561874801,1.029585650765054,0.9920501090198102,1.018546663416602,1.012703191458246",
    "22,1.02,1.018570612722768,1.025537301746956,0.9875604150241495,1.013862516598144,1.031262818747493,1.025533510270142,1.029287315065761,0.9914670000341481,1.017659271236473,1.012598655844716",
    "23,1.02,1.018873003721227,1.025765522675054,0.9879432794643023,1.014405176970968,1.031337139955975,1.025701661281636,1.029445560283458,0.991776070289318,1.018129764094604,1.012654097766168",
    "24,1.02,1.020063280702516,1.026663548903428,0.9894529299347244,1.016541130898494,1.031626390667237,1.026361970170326,1.030066518075306,0.9929938892766442,1.019980532117977,1.012871798785965",
    "25,1.02,1.021444078453059,1.027704673829624,0.9912096547607049,1.019019018781167,1.031955109984357,1.02712471881111,1.03078286366448,0.9944092447426414,1.022125236960042,1.013123246524867"
)

foreach ($line in $data) {
    $parts = $line.Split(",")
    $row = [int]$parts[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $colLetter = $cols[$i]
        $value = [double]$parts[$i + 1]
        $ws.Range($colLetter + $row).Value = $value
    }
}
